$d = $word.ActiveDocument

# --- 1. Create the three new character styles ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Underline = 1

# --- 2. Apply GaNStyle to every "Dates de la campanya 2022 ..." run (4 occurrences) ---

$datesText = "Dates de la campanya 2022 en què usem la  Constel·lació d'Hèrcules 13-22 de juny, 12-21 de juliol, 10-19 d'agost"

$rng = $d.Content
$rng.Start = 0
$found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
  $rng.Style = "GaNStyle"
  $rng.Collapse(0)
  $rng.End = $d.Content.End
  $found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- 3. Apply GaNParagraph to the "Esteu participant ..." run ---

$esteuText = "Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  Constel·lació d'Hèrcules a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn."

$rng2 = $d.Content
$rng2.Start = 0
$found2 = $rng2.Find.Execute($esteuText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
  $rng2.Style = "GaNParagraph"
}

# --- 4. Apply GaNLinks to the "Jenik Hollan, CzechGlobe (...)." run ---

$jenikText = "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng3 = $d.Content
$rng3.Start = 0
$found3 = $rng3.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
  $rng3.Style = "GaNLinks"
}
